$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row of price-tracking data for the latest scrape date.
# Force the cells to plain text so values are stored (and shared) as
# strings, matching the rest of the sheet rather than being
# auto-converted to dates/numbers.
$newRow = 38
$rng = $ws.Range("A$newRow:D$newRow")
$rng.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2026-02-07"
$ws.Cells.Item($newRow, 2).Value = "1149000"
$ws.Cells.Item($newRow, 3).Value = "8"
$ws.Cells.Item($newRow, 4).Value = "0"
